# Auto-generated: apply cryptos list price/volume refresh per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.641.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.62%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''1.883.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.27%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = '''1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''249.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +0.87%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = '''  +0.01%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = '''0.4753'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -0.28%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.2941'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +1.46%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''0.06538'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +0.19%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''21.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +0.65%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = '''0.07753'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.13%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''96.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +0.05%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''0.7386'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +0.23%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''1.878.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -0.06%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''5.246'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.33%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''274.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.52%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''30.611.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +0.51%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''13.19'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -3.29%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''0.000007538'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -0.91%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = '''  +0.02%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''2.127.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -0.02%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''5.349'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +2.07%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D24").Value = '''6.243'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.99%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''9.229'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -0.81%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''163.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.14%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''18.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.12%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''1.911'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.54%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = '''  -2.06%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''0.09733'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -2.33%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D32").Value = '''4.291'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -0.59%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''4.158'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +2.17%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''0.04867'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +1.93%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''1.126'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +0.16%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = '''0.7000'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.06%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = '''2.720'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +0.18%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''0.01909'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +2.11%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = '''  +2.27%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''6.302'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.44%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''74.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +5.85%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''2.031'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +4.74%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''0.4251'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.53%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''0.8410'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +0.37%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = '''  +0.03%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''102.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.00%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''9.369'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.05%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''7.055'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.45%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''35.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +0.27%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = '''917.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.96%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = '''0.05765'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +2.22%  '
$ws.Range("E51").Style = "Normal"
